# Apply the edits described by the commit diff:
#  1. Rename sheet "Employees" -> "Jul 2022"
#  2. Flip on the (already-default) "applyAlignment" flag for the two
#     plain/no-frills cell styles used by the header/footer filler cells,
#     without altering their visible alignment (WrapText is already False
#     for these cells, so re-asserting it forces the style engine to stamp
#     applyAlignment="1" on the xf record while keeping an empty <alignment/>,
#     matching the cleanup committed upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet rename -----------------------------------------------------
$ws.Name = "Jul 2022"

# --- Style clean-up -----------------------------------------------------
# Cells that shared cellXfs index 6 (numFmtId 0, no special formatting)
$plainCells = @("B1", "D1", "E1", "B2", "D2", "E2", "A97", "B97", "D97", "E97")
foreach ($addr in $plainCells) {
    $ws.Range($addr).WrapText = $false
}

# Cells that shared cellXfs index 7 (numFmtId 164, applyNumberFormat)
$dateCells = @("C1", "F1", "C2", "F2", "C97", "F97")
foreach ($addr in $dateCells) {
    $ws.Range($addr).WrapText = $false
}
